# "Generate Report for Archive"
#
# The localization status report is regenerated: the shared "Status" value
# moves from "Ready for handoff" to "In Translation" everywhere it is used
# (Overview!E2/F2 which mirror the zh-cn/de-de status, and the Status column
# on each per-locale sheet). Shortening that text causes Excel's column
# autosize to shrink the affected "Status" columns, so their widths are
# brought in line with the freshly generated report too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text (was "Ready for handoff").
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The shorter text narrows the autosized "Status" columns.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
